# edit.ps1
# Applies the cryptos.xlsx data-refresh update described in the commit
# "Updated cryptos list on Thu Mar 14 03:38:07 UTC 2024 with GitHub Actions"
#
# Every cell in this sheet (prices, % changes, coin name/link) is stored
# as TEXT, even when a price looks like a plain number (e.g. "1.00",
# "610.54"). Plain `Range.Value = "1.00"` lets Excel's COM layer coerce a
# numeric-looking string into a real number, which would change the cell
# type away from the text it must stay as. To avoid that, for any new
# value that parses as a number we briefly force the cell's number format
# to Text ("@") - the same effect as a user typing a leading apostrophe -
# assign the value, then clear the format back to the sheet's plain
# default so no stray formatting difference (e.g. a lingering "quote
# prefix" style) is left on the cell; the value itself stays text either
# way.
#
# Besides the per-cell value refresh, rows 28 and 29 swap contents: LEO
# moves up to row 28 and Filecoin moves down to row 29 (same four columns
# B/C/D/E updated either way), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "73.177.99"; ForceText = $False },
    @{ Cell = "E2"; Value = "  +1.64%  "; ForceText = $False },
    @{ Cell = "D3"; Value = "3.993.19"; ForceText = $False },
    @{ Cell = "D4"; Value = "1.00"; ForceText = $True },
    @{ Cell = "E4"; Value = "  +0.02%  "; ForceText = $False },
    @{ Cell = "D5"; Value = "610.54"; ForceText = $True },
    @{ Cell = "E5"; Value = "  +13.63%  "; ForceText = $False },
    @{ Cell = "D6"; Value = "166.18"; ForceText = $True },
    @{ Cell = "E6"; Value = "  +11.51%  "; ForceText = $False },
    @{ Cell = "D7"; Value = "0.683"; ForceText = $True },
    @{ Cell = "E7"; Value = "  -1.81%  "; ForceText = $False },
    @{ Cell = "E8"; Value = "  -0.07%  "; ForceText = $False },
    @{ Cell = "D9"; Value = "0.756"; ForceText = $True },
    @{ Cell = "E9"; Value = "  +0.79%  "; ForceText = $False },
    @{ Cell = "E10"; Value = "  +0.93%  "; ForceText = $False },
    @{ Cell = "D11"; Value = "56.90"; ForceText = $True },
    @{ Cell = "E11"; Value = "  +6.37%  "; ForceText = $False },
    @{ Cell = "D12"; Value = "0.0000338"; ForceText = $True },
    @{ Cell = "E12"; Value = "  +3.39%  "; ForceText = $False },
    @{ Cell = "D13"; Value = "11.10"; ForceText = $True },
    @{ Cell = "E13"; Value = "  +1.97%  "; ForceText = $False },
    @{ Cell = "D14"; Value = "4.630.53"; ForceText = $False },
    @{ Cell = "E14"; Value = "  -1.18%  "; ForceText = $False },
    @{ Cell = "D15"; Value = "3.989.33"; ForceText = $False },
    @{ Cell = "E15"; Value = "  -1.54%  "; ForceText = $False },
    @{ Cell = "E16"; Value = "  +4.71%  "; ForceText = $False },
    @{ Cell = "D17"; Value = "14.23"; ForceText = $True },
    @{ Cell = "E17"; Value = "  -0.40%  "; ForceText = $False },
    @{ Cell = "E18"; Value = "  -0.60%  "; ForceText = $False },
    @{ Cell = "D19"; Value = "73.052.71"; ForceText = $False },
    @{ Cell = "E19"; Value = "  +1.44%  "; ForceText = $False },
    @{ Cell = "E20"; Value = "  -0.22%  "; ForceText = $False },
    @{ Cell = "D21"; Value = "439.53"; ForceText = $True },
    @{ Cell = "E21"; Value = "  +0.63%  "; ForceText = $False },
    @{ Cell = "D22"; Value = "4.92"; ForceText = $True },
    @{ Cell = "E22"; Value = "  +15.53%  "; ForceText = $False },
    @{ Cell = "D23"; Value = "96.05"; ForceText = $True },
    @{ Cell = "E23"; Value = "  -2.08%  "; ForceText = $False },
    @{ Cell = "E24"; Value = "  -3.84%  "; ForceText = $False },
    @{ Cell = "D25"; Value = "14.22"; ForceText = $True },
    @{ Cell = "E25"; Value = "  -2.76%  "; ForceText = $False },
    @{ Cell = "E26"; Value = "  -6.56%  "; ForceText = $False },
    @{ Cell = "D27"; Value = "11.09"; ForceText = $True },
    @{ Cell = "E27"; Value = "  -1.52%  "; ForceText = $False },
    @{ Cell = "B28"; Value = "LEO"; ForceText = $False },
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; ForceText = $False },
    @{ Cell = "D28"; Value = "5.97"; ForceText = $True },
    @{ Cell = "E28"; Value = "  +0.40%  "; ForceText = $False },
    @{ Cell = "B29"; Value = "Filecoin"; ForceText = $False },
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; ForceText = $False },
    @{ Cell = "D29"; Value = "10.56"; ForceText = $True },
    @{ Cell = "E29"; Value = "  -1.44%  "; ForceText = $False },
    @{ Cell = "D30"; Value = "36.11"; ForceText = $True },
    @{ Cell = "E30"; Value = "  -2.54%  "; ForceText = $False },
    @{ Cell = "D31"; Value = "7.66"; ForceText = $True },
    @{ Cell = "E31"; Value = "  -7.92%  "; ForceText = $False },
    @{ Cell = "D32"; Value = "13.73"; ForceText = $True },
    @{ Cell = "E32"; Value = "  +1.27%  "; ForceText = $False },
    @{ Cell = "E33"; Value = "  -3.78%  "; ForceText = $False },
    @{ Cell = "D34"; Value = "72.11"; ForceText = $True },
    @{ Cell = "E34"; Value = "  +7.70%  "; ForceText = $False },
    @{ Cell = "E35"; Value = "  +19.78%  "; ForceText = $False },
    @{ Cell = "D36"; Value = "48.04"; ForceText = $True },
    @{ Cell = "E36"; Value = "  -4.21%  "; ForceText = $False },
    @{ Cell = "D37"; Value = "634.77"; ForceText = $True },
    @{ Cell = "E37"; Value = "  -7.13%  "; ForceText = $False },
    @{ Cell = "D38"; Value = "0.431"; ForceText = $True },
    @{ Cell = "E38"; Value = "  -6.31%  "; ForceText = $False },
    @{ Cell = "D39"; Value = "3.44"; ForceText = $True },
    @{ Cell = "E39"; Value = "  +1.41%  "; ForceText = $False },
    @{ Cell = "D40"; Value = "0.998"; ForceText = $True },
    @{ Cell = "E40"; Value = "  -0.31%  "; ForceText = $False },
    @{ Cell = "E41"; Value = "  -1.56%  "; ForceText = $False },
    @{ Cell = "D42"; Value = "11.05"; ForceText = $True },
    @{ Cell = "E42"; Value = "  -1.25%  "; ForceText = $False },
    @{ Cell = "E43"; Value = "  +0.18%  "; ForceText = $False },
    @{ Cell = "E44"; Value = "  -4.33%  "; ForceText = $False },
    @{ Cell = "D45"; Value = "0.0485"; ForceText = $True },
    @{ Cell = "E45"; Value = "  -1.55%  "; ForceText = $False },
    @{ Cell = "D46"; Value = "0.148"; ForceText = $True },
    @{ Cell = "E46"; Value = "  -1.44%  "; ForceText = $False },
    @{ Cell = "D47"; Value = "3.44"; ForceText = $True },
    @{ Cell = "E47"; Value = "  +4.07%  "; ForceText = $False },
    @{ Cell = "E48"; Value = "  -2.28%  "; ForceText = $False },
    @{ Cell = "D49"; Value = "2.87"; ForceText = $True },
    @{ Cell = "E49"; Value = "  +29.94%  "; ForceText = $False },
    @{ Cell = "D50"; Value = "2.877.60"; ForceText = $False },
    @{ Cell = "E50"; Value = "  +1.10%  "; ForceText = $False },
    @{ Cell = "D51"; Value = "3.03"; ForceText = $True },
    @{ Cell = "E51"; Value = "  -2.53%  "; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
